$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "GO"
$ws.Range("B8").Value = "TESTEGO"
$ws.Range("F8").Value = "T"
$ws.Range("H8").Value = "T - (T 02/11/25_12H) - GO"

$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "02/11/25"
